# Updated: st 04. 11. 2021
# Daily COVID stats refresh for OpenData_Slovakia_Covid_DailyStats (sheet1):
#  - A handful of historical AgTests (F) / AgPosit (G) values were revised.
#  - Two new daily rows (02. 11. 2021 / 03. 11. 2021) were appended.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Corrections to previously-reported AgTests / AgPosit figures ---
$ws.Range("F513").Value = 10582

$ws.Range("F556").Value = 12234

$ws.Range("F562").Value = 27098

$ws.Range("F567").Value = 23495

$ws.Range("F569").Value = 32442

$ws.Range("F572").Value = 33365

$ws.Range("F574").Value = 23425
$ws.Range("G574").Value = 355

$ws.Range("F575").Value = 26082

$ws.Range("F576").Value = 28966
$ws.Range("G576").Value = 433

$ws.Range("F580").Value = 28845
$ws.Range("G580").Value = 516

$ws.Range("F581").Value = 27061
$ws.Range("G581").Value = 480

$ws.Range("F582").Value = 25932
$ws.Range("G582").Value = 474

$ws.Range("F583").Value = 29277

$ws.Range("F586").Value = 33634

$ws.Range("F588").Value = 25355

$ws.Range("F589").Value = 25446
$ws.Range("G589").Value = 473

$ws.Range("F590").Value = 29299

$ws.Range("F591").Value = 14795
$ws.Range("G591").Value = 435

$ws.Range("F592").Value = 17945
$ws.Range("G592").Value = 632

$ws.Range("F594").Value = 29739

$ws.Range("F595").Value = 27236

$ws.Range("F596").Value = 29152

$ws.Range("F597").Value = 29422
$ws.Range("G597").Value = 954

$ws.Range("F598").Value = 15422
$ws.Range("G598").Value = 705

$ws.Range("F600").Value = 39721
$ws.Range("G600").Value = 1668

$ws.Range("F601").Value = 31485
$ws.Range("G601").Value = 1323

$ws.Range("F602").Value = 29464
$ws.Range("G602").Value = 1272

$ws.Range("F603").Value = 31316
$ws.Range("G603").Value = 1508

$ws.Range("F604").Value = 29185
$ws.Range("G604").Value = 1513

$ws.Range("F605").Value = 14415
$ws.Range("G605").Value = 1009

$ws.Range("F606").Value = 13701
$ws.Range("G606").Value = 1225

$ws.Range("F607").Value = 10554
$ws.Range("G607").Value = 932

# --- Append two new daily rows (608: 02.11.2021, 609: 03.11.2021) ---
$newRows = @(
    @(44502, 493277, 19150, 5361, 13112, 43251, 2765),
    @(44503, 499990, 21705, 6713, 13146, 26900, 1489)
)

$targetRow = 608
foreach ($rowValues in $newRows) {
    $ws.Cells.Item($targetRow, 1).Value = $rowValues[0]
    $ws.Cells.Item($targetRow, 2).Value = $rowValues[1]
    $ws.Cells.Item($targetRow, 3).Value = $rowValues[2]
    $ws.Cells.Item($targetRow, 4).Value = $rowValues[3]
    $ws.Cells.Item($targetRow, 5).Value = $rowValues[4]
    $ws.Cells.Item($targetRow, 6).Value = $rowValues[5]
    $ws.Cells.Item($targetRow, 7).Value = $rowValues[6]
    $targetRow++
}
